$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Trade #11 (Trade # 40 / row index 41 on "All Trades" & "MarketMaking",
# row 6 on "Strategy Status") closes at 2026-02-18 00:10:40 with ~0.000% P&L.
# ---------------------------------------------------------------------------

# --- Summary sheet -----------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.66   # Current Capital
$summary.Range("B4").Value = 0.76      # Total P&L $
$summary.Range("B5").Value = 0.39      # Total P&L %
$summary.Range("B6").Value = 39        # Total Trades
$summary.Range("B7").Value = 21        # Winning Trades
$summary.Range("B9").Value = 53.85     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.66
$status.Range("D6").Value = 10
$status.Range("E6").Value = -0.15
$status.Range("F6").Value = -0.34
$status.Range("G6").Value = 50

# --- All Trades sheet: close trade #40 (row 41) -------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G41").Value = 0.41
$allTrades.Range("H41").Value = "CLOSED"
$allTrades.Range("I41").Value = 7.8947
$allTrades.Range("J41").Value = 0.03
$allTrades.Range("K41").Value = 99.66
$allTrades.Range("L41").Value = "early_exit"
$allTrades.Range("M41").Value = 0.13

# --- All Trades sheet: append new trade #69 (row 70) ---------------------
$allTrades.Range("A70").Value = 69
$allTrades.Range("B70").Value = "'2026-02-18"
$allTrades.Range("C70").Value = "'00:10:35"
$allTrades.Range("D70").Value = "momentum"
$allTrades.Range("E70").Value = "UP"
$allTrades.Range("F70").Value = 0.38
$allTrades.Range("H70").Value = "OPEN"
$allTrades.Range("I70").Value = 0
$allTrades.Range("J70").Value = 0
$allTrades.Range("K70").Value = 100
$allTrades.Range("M70").Value = 0
$allTrades.Range("N70").Value = 0
$allTrades.Range("O70").Value = 0
$allTrades.Range("P70").Value = 0.9
$allTrades.Range("Q70").Value = "Upward momentum: 21.687% over 10 samples"

# --- momentum sheet: append new trade #69 (row 7) ------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A7").Value = 69
$momentum.Range("B7").Value = "'2026-02-18"
$momentum.Range("C7").Value = "'00:10:35"
$momentum.Range("D7").Value = "momentum"
$momentum.Range("E7").Value = "UP"
$momentum.Range("F7").Value = 0.38
$momentum.Range("H7").Value = "OPEN"
$momentum.Range("I7").Value = 0
$momentum.Range("J7").Value = 0
$momentum.Range("K7").Value = 100
$momentum.Range("L7").Value = 0
$momentum.Range("M7").Value = 0
$momentum.Range("N7").Value = 0.9
$momentum.Range("O7").Value = "Upward momentum: 21.687% over 10 samples"
$momentum.Range("Q7").Value = 0

# --- MarketMaking sheet: close trade #40 (row 12) -------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G12").Value = 0.41
$marketMaking.Range("H12").Value = "CLOSED"
$marketMaking.Range("I12").Value = 7.8947
$marketMaking.Range("J12").Value = 0.03
$marketMaking.Range("K12").Value = 99.66
$marketMaking.Range("P12").Value = "early_exit"
$marketMaking.Range("Q12").Value = 0.13
